$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Examtaker14301", "Automation14301", "examtakerautomation14301@gmail.com", "14301"),
    @("Examtaker33594", "Automation33594", "examtakerautomation33594@gmail.com", "33594"),
    @("Examtaker90878", "Automation90878", "examtakerautomation90878@gmail.com", "90878"),
    @("Examtaker55870", "Automation55870", "examtakerautomation55870@gmail.com", "55870"),
    @("Examtaker60396", "Automation60396", "examtakerautomation60396@gmail.com", "60396"),
    @("Examtaker75712", "Automation75712", "examtakerautomation75712@gmail.com", "75712"),
    @("Examtaker92863", "Automation92863", "examtakerautomation92863@gmail.com", "92863"),
    @("Examtaker54675", "Automation54675", "examtakerautomation54675@gmail.com", "54675"),
    @("Examtaker74526", "Automation74526", "examtakerautomation74526@gmail.com", "74526"),
    @("Examtaker65907", "Automation65907", "examtakerautomation65907@gmail.com", "65907"),
    @("Examtaker40646", "Automation40646", "examtakerautomation40646@gmail.com", "40646"),
    @("Examtaker34537", "Automation34537", "examtakerautomation34537@gmail.com", "34537"),
    @("Examtaker41965", "Automation41965", "examtakerautomation41965@gmail.com", "41965"),
    @("Examtaker95673", "Automation95673", "examtakerautomation95673@gmail.com", "95673"),
    @("Examtaker43730", "Automation43730", "examtakerautomation43730@gmail.com", "43730"),
    @("Examtaker99912", "Automation99912", "examtakerautomation99912@gmail.com", "99912"),
    @("Examtaker31001", "Automation31001", "examtakerautomation31001@gmail.com", "31001"),
    @("Examtaker38723", "Automation38723", "examtakerautomation38723@gmail.com", "38723"),
    @("Examtaker30939", "Automation30939", "examtakerautomation30939@gmail.com", "30939"),
    @("Examtaker97646", "Automation97646", "examtakerautomation97646@gmail.com", "97646"),
    @("Examtaker11135", "Automation11135", "examtakerautomation11135@gmail.com", "11135"),
    @("Examtaker44337", "Automation44337", "examtakerautomation44337@gmail.com", "44337"),
    @("Examtaker98912", "Automation98912", "examtakerautomation98912@gmail.com", "98912"),
    @("Examtaker28571", "Automation28571", "examtakerautomation28571@gmail.com", "28571"),
    @("Examtaker28338", "Automation28338", "examtakerautomation28338@gmail.com", "28338"),
    @("Examtaker65991", "Automation65991", "examtakerautomation65991@gmail.com", "65991"),
    @("Examtaker33199", "Automation33199", "examtakerautomation33199@gmail.com", "33199"),
    @("Examtaker79128", "Automation79128", "examtakerautomation79128@gmail.com", "79128"),
    @("Examtaker74825", "Automation74825", "examtakerautomation74825@gmail.com", "74825"),
    @("Examtaker84552", "Automation84552", "examtakerautomation84552@gmail.com", "84552")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $cellA = $ws.Cells.Item($r, 1)
    $cellB = $ws.Cells.Item($r, 2)
    $cellC = $ws.Cells.Item($r, 3)
    $cellD = $ws.Cells.Item($r, 4)
    $cellA.Value = $data[$i][0]
    $cellB.Value = $data[$i][1]
    $cellC.Value = $data[$i][2]
    $cellD.NumberFormat = "@"
    $cellD.Value = $data[$i][3]
    $cellD.Style = "Normal"
}
